# Updated cryptos list values (price / volume(1h) columns, plus the
# Monero/Cosmos row swap) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '30.451.53' }
    @{ Cell = 'E2'; Value = '  +1.02%  ' }
    @{ Cell = 'D3'; Value = '1.853.09' }
    @{ Cell = 'E3'; Value = '  +1.23%  ' }
    @{ Cell = 'E4'; Value = '  -0.02%  ' }
    @{ Cell = 'D5'; Value = '233.44' }
    @{ Cell = 'E5'; Value = '  +1.02%  ' }
    @{ Cell = 'E6'; Value = '  -0.01%  ' }
    @{ Cell = 'D7'; Value = '0.4757' }
    @{ Cell = 'E7'; Value = '  +2.27%  ' }
    @{ Cell = 'D8'; Value = '0.2758' }
    @{ Cell = 'E8'; Value = '  +2.72%  ' }
    @{ Cell = 'D9'; Value = '0.06348' }
    @{ Cell = 'E9'; Value = '  +1.40%  ' }
    @{ Cell = 'D10'; Value = '17.98' }
    @{ Cell = 'E10'; Value = '  +12.43%  ' }
    @{ Cell = 'D11'; Value = '1.893.86' }
    @{ Cell = 'E11'; Value = '  +3.08%  ' }
    @{ Cell = 'E12'; Value = '  +1.11%  ' }
    @{ Cell = 'D13'; Value = '4.971' }
    @{ Cell = 'E14'; Value = '  +1.98%  ' }
    @{ Cell = 'D15'; Value = '0.6248' }
    @{ Cell = 'E15'; Value = '  +1.20%  ' }
    @{ Cell = 'D16'; Value = '30.416.93' }
    @{ Cell = 'E16'; Value = '  +1.16%  ' }
    @{ Cell = 'D17'; Value = '244.81' }
    @{ Cell = 'E19'; Value = '  +2.68%  ' }
    @{ Cell = 'D20'; Value = '0.000007354' }
    @{ Cell = 'E20'; Value = '  +1.18%  ' }
    @{ Cell = 'D21'; Value = '0.9990' }
    @{ Cell = 'E21'; Value = '  -0.13%  ' }
    @{ Cell = 'D22'; Value = '4.931' }
    @{ Cell = 'E22'; Value = '  +1.99%  ' }
    @{ Cell = 'D23'; Value = '5.919' }
    @{ Cell = 'E23'; Value = '  +1.11%  ' }
    @{ Cell = 'B24'; Value = 'Cosmos' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D24'; Value = '9.094' }
    @{ Cell = 'E24'; Value = '  -0.17%  ' }
    @{ Cell = 'B25'; Value = 'Monero' }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D25'; Value = '164.32' }
    @{ Cell = 'E25'; Value = '  -0.02%  ' }
    @{ Cell = 'E26'; Value = '  +2.23%  ' }
    @{ Cell = 'D27'; Value = '1.879' }
    @{ Cell = 'E27'; Value = '  +2.03%  ' }
    @{ Cell = 'D28'; Value = '0.1028' }
    @{ Cell = 'E28'; Value = '  +2.07%  ' }
    @{ Cell = 'E29'; Value = '  -1.39%  ' }
    @{ Cell = 'D30'; Value = '4.045' }
    @{ Cell = 'E30'; Value = '  +0.19%  ' }
    @{ Cell = 'D31'; Value = '3.842' }
    @{ Cell = 'E31'; Value = '  +2.27%  ' }
    @{ Cell = 'E32'; Value = '  +1.37%  ' }
    @{ Cell = 'E33'; Value = '  +0.70%  ' }
    @{ Cell = 'D34'; Value = '0.6992' }
    @{ Cell = 'E34'; Value = '  -0.07%  ' }
    @{ Cell = 'D35'; Value = '2.700' }
    @{ Cell = 'E35'; Value = '  +0.47%  ' }
    @{ Cell = 'D36'; Value = '0.01898' }
    @{ Cell = 'E36'; Value = '  +4.98%  ' }
    @{ Cell = 'D37'; Value = '2.679' }
    @{ Cell = 'E37'; Value = '  +2.67%  ' }
    @{ Cell = 'D38'; Value = '0.8791' }
    @{ Cell = 'E38'; Value = '  -1.28%  ' }
    @{ Cell = 'D40'; Value = '106.66' }
    @{ Cell = 'E40'; Value = '  +3.90%  ' }
    @{ Cell = 'D42'; Value = '0.4074' }
    @{ Cell = 'E42'; Value = '  +2.11%  ' }
    @{ Cell = 'D43'; Value = '5.502' }
    @{ Cell = 'E43'; Value = '  +0.73%  ' }
    @{ Cell = 'D44'; Value = '7.182' }
    @{ Cell = 'E44'; Value = '  +3.47%  ' }
    @{ Cell = 'D45'; Value = '63.49' }
    @{ Cell = 'E45'; Value = '  +6.69%  ' }
    @{ Cell = 'E46'; Value = '  +1.19%  ' }
    @{ Cell = 'D47'; Value = '33.98' }
    @{ Cell = 'E47'; Value = '  +4.35%  ' }
    @{ Cell = 'D48'; Value = '8.606' }
    @{ Cell = 'E48'; Value = '  +1.88%  ' }
    @{ Cell = 'E49'; Value = '  -0.26%  ' }
    @{ Cell = 'D50'; Value = '1.350' }
    @{ Cell = 'E50'; Value = '  -0.51%  ' }
    @{ Cell = 'D51'; Value = '0.3695' }
    @{ Cell = 'E51'; Value = '  +2.29%  ' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    # Preserve the existing style while forcing the assignment to be
    # stored as text -- several values look numeric (e.g. "233.44",
    # "0.9990", "1.350") and Excel would otherwise coerce them to
    # numbers and silently drop significant trailing zeros / dots.
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = $origStyle
}
